$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.013.94"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.270.47"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.32"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.96"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.20"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.24"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.615.12"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.276.75"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.56"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.801"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.939.56"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  +10.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0916"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.58"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.86"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.87"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.69"
$ws.Range("E27").Value = "  +11.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.59"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.67"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.63"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.54"
$ws.Range("E32").Value = "  -6.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0793"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.56"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.94"
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.103"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.94"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0311"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.63"
$ws.Range("E42").Value = "  -10.44%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("E44").Value = "  +12.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.760.53"
$ws.Range("E45").Value = "  -3.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.193"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "76.62"
$ws.Range("E47").Value = "  -4.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "69.73"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.62"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.72"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.70"
$ws.Range("E51").Value = "  -3.24%  "
